$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E keep their existing text formatting, so
# numeric-looking strings (e.g. "96.99") are not silently converted
# into actual numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "43.736.26"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.309.02"
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "96.99"
$ws.Range("E5").Value = "  +3.72%  "
$ws.Range("D6").Value = "272.72"
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").Value = "45.20"
$ws.Range("E10").Value = "  -1.87%  "
$ws.Range("D11").Value = "0.0953"
$ws.Range("E11").Value = "  -2.29%  "
$ws.Range("E12").Value = "  -3.57%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "2.646.84"
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("D15").Value = "15.51"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("D16").Value = "0.872"
$ws.Range("E16").Value = "  +7.39%  "
$ws.Range("D17").Value = "2.319.31"
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("D18").Value = "43.706.69"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("E20").Value = "  +4.08%  "
$ws.Range("D21").Value = "73.50"
$ws.Range("E21").Value = "  +3.61%  "
$ws.Range("D22").Value = "239.44"
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").Value = "2.28"
$ws.Range("E23").Value = "  -2.69%  "
$ws.Range("D24").Value = "9.44"
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").Value = "11.35"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "3.49"
$ws.Range("E28").Value = "  -2.25%  "
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "38.21"
$ws.Range("E30").Value = "  -7.19%  "
$ws.Range("E31").Value = "  +6.75%  "
$ws.Range("D32").Value = "175.21"
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("D33").Value = "0.0915"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").Value = "5.48"
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("E36").Value = "  +3.54%  "
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("D38").Value = "4.45"
$ws.Range("E38").Value = "  +2.48%  "
$ws.Range("D39").Value = "3.37"
$ws.Range("E39").Value = "  -6.20%  "
$ws.Range("D40").Value = "0.245"
$ws.Range("E40").Value = "  +8.10%  "
$ws.Range("D41").Value = "2.40"
$ws.Range("E41").Value = "  +10.25%  "
$ws.Range("D42").Value = "1.41"
$ws.Range("E42").Value = "  +21.68%  "
$ws.Range("E43").Value = "  -5.13%  "
$ws.Range("D44").Value = "62.97"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("D45").Value = "9.18"
$ws.Range("E45").Value = "  +9.31%  "
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("E47").Value = "  +4.02%  "
$ws.Range("D48").Value = "100.42"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("E50").Value = "  +14.43%  "
$ws.Range("D51").Value = "2.535.35"
$ws.Range("E51").Value = "  +2.79%  "
